$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 470, shifting existing rows 470:562 down to 471:563
$ws.Rows.Item(470).Insert()

# Populate the newly inserted row 470 with the new data entry
$ws.Cells.Item(470, 1).Value = 10
$ws.Cells.Item(470, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(470, 3).Value = "La Araucanía"
$ws.Cells.Item(470, 4).Value = 45258
$ws.Cells.Item(470, 5).Value = 9
$ws.Cells.Item(470, 6).Value = 100114013
$ws.Cells.Item(470, 7).Value = "Zanahoria"
$ws.Cells.Item(470, 8).Value = "Sin especificar"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 155
$ws.Cells.Item(470, 11).Value = 5000
$ws.Cells.Item(470, 12).Value = 5000
$ws.Cells.Item(470, 13).Value = 5000
$ws.Cells.Item(470, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(470, 15).Value = "Región del Maule"
$ws.Cells.Item(470, 16).Value = 250
$ws.Cells.Item(470, 17).Value = 20
$ws.Cells.Item(470, 18).Value = "Hortaliza"
